$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Inflação argentina volta a subir após 5 meses e bate 271,5% em 1 ano"
$ws.Range("B2").Value = " "
$ws.Range("C2").Value = "https://g1.globo.com/economia/noticia/2024/07/12/inflacao-argentina.ghtml"

$ws.Range("A3").Value = "PF marca depoimento de Ramagem por espionagem no governo Bolsonaro"
$ws.Range("B3").Value = " "
$ws.Range("C3").Value = "https://g1.globo.com/politica/noticia/2024/07/12/pf-marca-para-quarta-feira-depoimento-de-ramagem-no-caso-da-abin-paralela.ghtml"

$ws.Range("A4").Value = "ANÁLISE: entrevista de Biden foi teste de cognição constrangedor"
$ws.Range("B4").Value = "Presidente dos EUA teve desempenho superior ao do debate, mas gafes roubaram a cena."
$ws.Range("C4").Value = "https://g1.globo.com/mundo/blog/sandra-cohen/post/2024/07/12/biden-enfrenta-entrevista-como-um-constrangedor-teste-de-cognicao.ghtml"

$ws.Range("A5").Value = "Governo avalia que 'Jogo do Tigrinho' pode ser liberado no Brasil"
$ws.Range("B5").Value = " "
$ws.Range("C5").Value = "https://g1.globo.com/tecnologia/noticia/2024/07/12/governo-fortune-tiger-sites-exterior.ghtml"

$ws.Range("A6").Value = "Golpe do 'PIX errado': saiba como não ser enganado"
$ws.Range("B6").Value = " "
$ws.Range("C6").Value = "https://g1.globo.com/economia/noticia/2024/07/12/golpe-do-pix-errado-saiba-como-os-criminosos-agem-e-como-nao-ser-enganado.ghtml"

$ws.Range("A7").Value = "Polícia indicia 6 por morte de empresário com doce envenenado"
$ws.Range("B7").Value = "Namorada e mulher que se apresenta como cigana vão responder por homicídio."
$ws.Range("C7").Value = "https://g1.globo.com/rj/rio-de-janeiro/noticia/2024/07/12/policia-morte-de-empresario-envenenado-com-brigadeirao-no-rio.ghtml"

$ws.Range("A8").Value = "GPS ajuda polícia dos EUA a prender suspeito de matar brasileira"
$ws.Range("B8").Value = "Corpo de Suzan Ferreira foi achado às margens de estrada em Michigan."
$ws.Range("C8").Value = "https://g1.globo.com/mg/minas-gerais/noticia/2024/07/12/policia-americana-prende-suspeito-de-envolvimento-na-morte-de-brasileira-corpo-foi-encontrado-em-rodovia.ghtml"

$ws.Range("A9").Value = "Nego Di é suspeito de lavagem de R$ 2 milhões com rifas; mulher é presa"
$ws.Range("B9").Value = " "
$ws.Range("C9").Value = "https://g1.globo.com/rs/rio-grande-do-sul/noticia/2024/07/12/nego-di-e-alvo-de-operacao-que-investiga-suspeita-de-lavagem-de-r-2-milhoes-com-rifas-virtuais.ghtml"
